# "Changes of 5th May 2022"
# Replace the package tracking numbers in column C (rows 2-22) with a new
# batch of numbers. Rows whose column D mirrors column C (5,6,7,13-17) get
# the same new value written into D as well.
#
# Values are entered with a leading apostrophe so Excel stores them as text
# (these are long digit strings that must stay exactly as typed, not be
# coerced into floating point / scientific notation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTrackNums = @{
    2  = "320018471602"
    3  = "320018471613"
    4  = "320018471646"
    5  = "320018471668"
    6  = "320018471705"
    7  = "320018471727"
    8  = "320018471750"
    9  = "320018471771"
    10 = "320018471808"
    11 = "320018471820"
    12 = "320018471863"
    13 = "320018471885"
    14 = "320018471911"
    15 = "320018471933"
    16 = "320018471966"
    17 = "320018471988"
    18 = "320018472024"
    19 = "320018472046"
    20 = "320018472079"
    21 = "320018472090"
    22 = "320018472127"
}

# Rows where column D duplicates column C's tracking number.
$rowsWithD = @(5, 6, 7, 13, 14, 15, 16, 17)

foreach ($row in 2..22) {
    $value = "'" + $newTrackNums[$row]

    $ws.Range("C$row").Value = $value

    if ($rowsWithD -contains $row) {
        $ws.Range("D$row").Value = $value
    }
}
